$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = 1
$ws.Range("C70").Value = "2024-06-16 09:14:15"
$ws.Range("D70").Value = 200
$ws.Range("E70").Value = 8

# Row 71
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = 2
$ws.Range("C71").Value = "2024-06-16 09:14:16"
$ws.Range("D71").Value = 200
$ws.Range("E71").Value = 0
